$d = $word.ActiveDocument

$replacements = @(
    @("24×18=432", "49×38=1862"),
    @("26×65=1690", "59×12=708"),
    @("36×13=468", "48×47=2256"),
    @("51×56=2856", "62×78=4836"),
    @("59×67=3953", "62×49=3038"),
    @("72×58=4176", "68×33=2244"),
    @("64×50=3200", "87×47=4089"),
    @("14×38=532", "67×94=6298"),
    @("87×55=4785", "62×51=3162"),
    @("60×99=5940", "84×11=924"),
    @("65×62=4030", "17×21=357"),
    @("43×47=2021", "73×79=5767"),
    @("47×71=3337", "42×18=756"),
    @("19×84=1596", "27×62=1674"),
    @("73×28=2044", "90×49=4410"),
    @("66×11=726", "47×72=3384"),
    @("78×84=6552", "57×79=4503"),
    @("24×27=648", "61×16=976"),
    @("81×23=1863", "87×24=2088"),
    @("96×68=6528", "24×15=360"),
    @("19×83=1577", "11×64=704"),
    @("25×55=1375", "11×27=297"),
    @("36×70=2520", "37×70=2590"),
    @("76×92=6992", "80×16=1280"),
    @("14×97=1358", "77×16=1232")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
